# Append a new blank line plus two new paragraphs (date/name line and the
# bio paragraph) to the end of the document, after David Singletary's
# "... in our Data Science program." paragraph -- mirroring a second
# student ("Lillian MacKenzie") adding her icebreaker entry per the
# assignment instructions ("add yours to the end of the file").

$d = $word.ActiveDocument

function Get-EndRange {
    $r = $d.Content
    $r.Collapse(0)
    return $r
}

# 1) Blank paragraph separating the two entries.
$r = Get-EndRange
$r.InsertParagraphAfter()

# 2) "9/05/2023 Lillian MacKenzie" paragraph.
$r = Get-EndRange
$r.InsertParagraphAfter()
$r = Get-EndRange
$r.InsertAfter("9/05/2023 Lillian MacKenzie")

# 3) Bio paragraph.
$r = Get-EndRange
$r.InsertParagraphAfter()
$r = Get-EndRange
$r.InsertAfter("Hello, my name is Lillian MacKenzie, but I prefer to go by Luke. I" + [char]8217 + "ve lived in Middleburg Florida my whole life. I love games and classic sci-fi movies, especially Star Wars and Lord of the Rings. One day I want to work in automation.  ")

Write-Output "Appended Lillian MacKenzie entry."
